$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark (an internal marker of the last
#    edit position; Word normally drops/relocates it as a side effect of
#    saving and it is not meaningful document content).
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2. Locate the paragraph that ends in "... button_debounce component."
#    and append a new sentence (as its own run) telling students which
#    template file to use for their design.
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("button_debounce component.")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pEnd = $target.Range.End
    # Collapse to a zero-length point just before the paragraph mark so the
    # inserted text lands inside the paragraph, as a new trailing run.
    $insertPoint = $d.Range($pEnd - 1, $pEnd - 1)
    $insertPoint.InsertAfter(" For your design, use HW7_design_template.pptx.")
}
